$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D and E columns keep their text formatting so numeric-looking
# strings (e.g. "35.30", "0.0790") are not coerced into numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "47.801.65"
$ws.Range("E2").Value = "  +1.29%  "

$ws.Range("D3").Value = "2.494.50"
$ws.Range("E3").Value = "  +0.16%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").Value = "323.12"
$ws.Range("E5").Value = "  +0.16%  "

$ws.Range("D6").Value = "108.94"
$ws.Range("E6").Value = "  +0.95%  "

$ws.Range("E7").Value = "  -0.41%  "

$ws.Range("E8").Value = "  +0.09%  "

$ws.Range("D9").Value = "0.551"
$ws.Range("E9").Value = "  +1.57%  "

$ws.Range("D10").Value = "40.69"
$ws.Range("E10").Value = "  +6.64%  "

$ws.Range("E11").Value = "  +0.12%  "

$ws.Range("E12").Value = "  +0.63%  "

$ws.Range("D13").Value = "18.64"
$ws.Range("E13").Value = "  +1.29%  "

$ws.Range("E14").Value = "  +0.20%  "

$ws.Range("D15").Value = "2.882.74"
$ws.Range("E15").Value = "  -0.06%  "

$ws.Range("D16").Value = "2.488.31"
$ws.Range("E16").Value = "  -0.49%  "

$ws.Range("D17").Value = "0.851"
$ws.Range("E17").Value = "  -0.34%  "

$ws.Range("D18").Value = "47.697.25"
$ws.Range("E18").Value = "  +1.21%  "

$ws.Range("D19").Value = "13.22"
$ws.Range("E19").Value = "  +2.53%  "

$ws.Range("D20").Value = "6.63"
$ws.Range("E20").Value = "  -0.83%  "

$ws.Range("E21").Value = "  +0.17%  "

$ws.Range("E22").Value = "  +13.82%  "

$ws.Range("E23").Value = "  +0.28%  "

$ws.Range("D24").Value = "247.39"
$ws.Range("E24").Value = "  -1.29%  "

$ws.Range("E25").Value = "  -1.77%  "

$ws.Range("E26").Value = "  +0.13%  "

$ws.Range("D27").Value = "25.86"
$ws.Range("E27").Value = "  -1.11%  "

$ws.Range("D28").Value = "9.99"
$ws.Range("E28").Value = "  -0.31%  "

$ws.Range("E29").Value = "  -0.76%  "

$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D30").Value = "35.30"
$ws.Range("E30").Value = "  +1.10%  "

$ws.Range("B31").Value = "Kaspa"
$ws.Range("C31").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D31").Value = "0.139"
$ws.Range("E31").Value = "  +0.40%  "

$ws.Range("D32").Value = "49.86"
$ws.Range("E32").Value = "  +1.02%  "

$ws.Range("D33").Value = "19.91"
$ws.Range("E33").Value = "  +1.68%  "

$ws.Range("E34").Value = "  -2.23%  "

$ws.Range("D35").Value = "0.0790"
$ws.Range("E35").Value = "  +0.22%  "

$ws.Range("E36").Value = "  +0.25%  "

$ws.Range("E37").Value = "  -1.02%  "

$ws.Range("E38").Value = "  +0.13%  "

$ws.Range("E39").Value = "  -1.22%  "

$ws.Range("E40").Value = "  -0.15%  "

$ws.Range("D41").Value = "22.43"
$ws.Range("E41").Value = "  +5.76%  "

$ws.Range("E42").Value = "  -1.00%  "

$ws.Range("D43").Value = "119.02"
$ws.Range("E43").Value = "  -2.27%  "

$ws.Range("E44").Value = "  +0.05%  "

$ws.Range("D45").Value = "2.002.73"
$ws.Range("E45").Value = "  +1.97%  "

$ws.Range("D46").Value = "3.05"
$ws.Range("E46").Value = "  +1.49%  "

$ws.Range("E47").Value = "  -3.48%  "

$ws.Range("E48").Value = "  +1.18%  "

$ws.Range("D49").Value = "9.02"
$ws.Range("E49").Value = "  -0.30%  "

$ws.Range("D50").Value = "5.14"
$ws.Range("E50").Value = "  -2.57%  "

$ws.Range("D51").Value = "56.78"
